$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 121
$ws.Range("B3").Value = 650
$ws.Range("B4").Value = 50
$ws.Range("B5").Value = 400
$ws.Range("B6").Value = 850
